# Automatic update 2025-07-07 14:30:08
#
# Inserts 5 new clients (in alphabetical order) for advisor "RIOS CARRION
# ANGEL BENIGNO" into the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets:
#   CONSTANTE CAMACHO ARIANA ELIZABETH
#   CULMA OVIEDO NINI JOHANA
#   LUNA PAZMIÑO MYRIAM DEL ROCIO
#   MERIZALDE PEREIRA KAREN ELIZABETH
#   MUNDIACABADOS CIA. LTDA.
# and updates the totals row accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A..R, product-group breakdown)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert blank rows at the original (pre-edit) positions, working from the
# bottom up so earlier row numbers stay valid while we work:
#  - 3 blank rows before original row 13 (MUNDO-CERAMICO CIA.LTDA.)
#  - 1 blank row before original row 9  (DDH S.A.S.)
#  - 1 blank row before original row 8  (CONZA VEGA FRANCO BLADYMIR)
$ws1.Rows.Item(13).Resize(3).Insert()
$ws1.Rows.Item(9).Insert()
$ws1.Rows.Item(8).Insert()

# New client rows (all product-group quantities are 0 for these new clients)
$newRows1 = @(
    @{ Row = 8;  Name = "CONSTANTE CAMACHO ARIANA ELIZABETH" },
    @{ Row = 10; Name = "CULMA OVIEDO NINI JOHANA" },
    @{ Row = 15; Name = "LUNA PAZMIÑO MYRIAM DEL ROCIO" },
    @{ Row = 16; Name = "MERIZALDE PEREIRA KAREN ELIZABETH" },
    @{ Row = 17; Name = "MUNDIACABADOS CIA. LTDA." }
)

foreach ($entry in $newRows1) {
    $r = $entry.Row
    $ws1.Range("A$r").Value = "RIOS CARRION ANGEL BENIGNO"
    $ws1.Range("B$r").Value = $entry.Name
    for ($col = 3; $col -le 18; $col++) {
        $ws1.Cells.Item($r, $col).Value = 0
    }
}

# The summary row (was row 19, "N de 17") is now row 24 and must reflect the
# new client count (22 instead of 17).
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(24, $col)
    $cell.Value = $cell.Text -replace "de 17", "de 22"
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A..G, month-by-month sales)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(13).Resize(3).Insert()
$ws2.Rows.Item(9).Insert()
$ws2.Rows.Item(8).Insert()

$newRows2 = @(
    @{ Row = 8;  Name = "CONSTANTE CAMACHO ARIANA ELIZABETH";  C = 0;       D = 0; E = 2261.64;  F = 0; G = 0 },
    @{ Row = 10; Name = "CULMA OVIEDO NINI JOHANA";            C = 0;       D = 0; E = 851.43;    F = 0; G = 0 },
    @{ Row = 15; Name = "LUNA PAZMIÑO MYRIAM DEL ROCIO";       C = 0;       D = 0; E = 0;         F = 0; G = 0 },
    @{ Row = 16; Name = "MERIZALDE PEREIRA KAREN ELIZABETH";   C = 0;       D = 0; E = 0;         F = 0; G = 0 },
    @{ Row = 17; Name = "MUNDIACABADOS CIA. LTDA.";            C = -545.18; D = 0; E = 0;         F = 0; G = 0 }
)

foreach ($entry in $newRows2) {
    $r = $entry.Row
    $ws2.Range("A$r").Value = "RIOS CARRION ANGEL BENIGNO"
    $ws2.Range("B$r").Value = $entry.Name
    $ws2.Range("C$r").Value = $entry.C
    $ws2.Range("D$r").Value = $entry.D
    $ws2.Range("E$r").Value = $entry.E
    $ws2.Range("F$r").Value = $entry.F
    $ws2.Range("G$r").Value = $entry.G
}

# Recompute the totals row (was row 19, now row 24) across all 22 clients.
# NOTE: reading back through ".Value" surfaces a reflection placeholder in
# this host; ".Value2" returns the real numeric payload, so use that for
# the read side of this round-trip (writes still use ".Value").
$totalRow = 24
foreach ($col in @("C", "D", "E", "F", "G")) {
    $sum = 0
    for ($r = 2; $r -le 23; $r++) {
        $sum += [double]$ws2.Range("$col$r").Value2
    }
    $ws2.Range("$col$totalRow").Value = $sum
}
